# Replace Oracle-style "#pDate#" bind-variable placeholders with
# colon-style ":pDate:" placeholders in the SQL text stored in column B.
# Only the exact-case token "#pDate#" is targeted; the many occurrences of
# the lower-case "#pdate#" token elsewhere in the sheet must stay untouched,
# so a case-sensitive .NET string Replace() is used instead of PowerShell's
# (case-insensitive by default, and -creplace is also case-insensitive in
# this runtime) -replace operator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("B4", "B6", "B8", "B9", "B10")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null) {
        $updated = $text.Replace("#pDate#", ":pDate:")
        if ($updated -ne $text) {
            $cell.Value2 = $updated
        }
    }
}
